$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 113, shifting existing rows 113:206 down to 114:207.
$ws.Rows("113:113").Insert()

# Populate the newly inserted row 113 with the new record's data.
$ws.Range("A113").Value = 3
$ws.Range("B113").Value = "Femacal de La Calera"
$ws.Range("C113").Value = "Coquimbo"
$ws.Range("D113").Value = 44907
$ws.Range("E113").Value = 5
$ws.Range("F113").Value = 100112030
$ws.Range("G113").Value = "Poroto granado"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 73
$ws.Range("K113").Value = 31000
$ws.Range("L113").Value = 32080
$ws.Range("M113").Value = 31562
$ws.Range("N113").Value = "$/malla 25 kilos"
$ws.Range("O113").Value = "Provincia de Limarí"
$ws.Range("P113").Value = 1262
$ws.Range("Q113").Value = 25
$ws.Range("R113").Value = "Hortaliza"
